$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param(
        [string]$CellRef,
        [string]$CellValue
    )
    $rng = $ws.Range($CellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $CellValue
    $rng.Style = "Normal"
}

# Row 2
Set-TextCell 'D2' '60.502.09'
Set-TextCell 'E2' '  +1.69%  '

# Row 3
Set-TextCell 'D3' '2.604.90'
Set-TextCell 'E3' '  +0.67%  '

# Row 4
Set-TextCell 'E4' '  -0.02%  '

# Row 5
Set-TextCell 'D5' '573.08'
Set-TextCell 'E5' '  +1.96%  '

# Row 6
Set-TextCell 'D6' '142.83'
Set-TextCell 'E6' '  -0.16%  '

# Row 7
Set-TextCell 'E7' '  -0.24%  '

# Row 8
Set-TextCell 'E8' '  +0.53%  '

# Row 9
Set-TextCell 'D9' '2.629.24'
Set-TextCell 'E9' '  +1.27%  '

# Row 10
Set-TextCell 'D10' '6.49'
Set-TextCell 'E10' '  -2.81%  '

# Row 11
Set-TextCell 'D11' '0.106'
Set-TextCell 'E11' '  +1.04%  '

# Row 12
Set-TextCell 'E12' '  -4.48%  '

# Row 13
Set-TextCell 'E13' '  +2.11%  '

# Row 14
Set-TextCell 'D14' '3.069.11'
Set-TextCell 'E14' '  +0.89%  '

# Row 15
Set-TextCell 'D15' '60.507.51'
Set-TextCell 'E15' '  +1.83%  '

# Row 16
Set-TextCell 'D16' '23.23'
Set-TextCell 'E16' '  -0.42%  '

# Row 17
Set-TextCell 'E17' '  +2.59%  '

# Row 18
Set-TextCell 'D18' '2.637.67'
Set-TextCell 'E18' '  +1.44%  '

# Row 19
Set-TextCell 'D19' '11.33'

# Row 20
Set-TextCell 'E20' '  +1.79%  '

# Row 21
Set-TextCell 'D21' '347.31'
Set-TextCell 'E21' '  +2.48%  '

# Row 22
Set-TextCell 'E22' '  +6.59%  '

# Row 23
Set-TextCell 'D23' '0.997'
Set-TextCell 'E23' '  -0.34%  '

# Row 24
Set-TextCell 'D24' '0.526'
Set-TextCell 'E24' '  +11.26%  '

# Row 25
Set-TextCell 'D25' '63.23'
Set-TextCell 'E25' '  -0.25%  '

# Row 26
Set-TextCell 'E26' '  -0.04%  '

# Row 27
Set-TextCell 'E27' '  -1.51%  '

# Row 28
Set-TextCell 'D28' '7.75'
Set-TextCell 'E28' '  +3.70%  '

# Row 29
Set-TextCell 'D29' '0.0₃0789'
Set-TextCell 'E29' '  +1.13%  '

# Row 30
Set-TextCell 'D30' '1.85'
Set-TextCell 'E30' '  +10.00%  '

# Row 31
Set-TextCell 'E31' '  +2.23%  '

# Row 32
Set-TextCell 'E32' '  -0.07%  '

# Row 33
Set-TextCell 'D33' '161.50'
Set-TextCell 'E33' '  +1.94%  '

# Row 34
Set-TextCell 'D34' '19.52'
Set-TextCell 'E34' '  +2.26%  '

# Row 35
Set-TextCell 'D35' '4.23'
Set-TextCell 'E35' '  +4.08%  '

# Row 36
Set-TextCell 'D36' '0.975'
Set-TextCell 'E36' '  +8.38%  '

# Row 37
Set-TextCell 'E37' '  +4.23%  '

# Row 38
Set-TextCell 'E38' '  +7.53%  '

# Row 39
Set-TextCell 'E39' '  +1.03%  '

# Row 40
Set-TextCell 'D40' '3.85'
Set-TextCell 'E40' '  +4.39%  '

# Row 41
Set-TextCell 'D41' '0.841'
Set-TextCell 'E41' '  -3.44%  '

# Row 42
Set-TextCell 'D42' '295.24'
Set-TextCell 'E42' '  +0.05%  '

# Row 43
Set-TextCell 'D43' '137.29'
Set-TextCell 'E43' '  -1.86%  '

# Row 44
Set-TextCell 'E44' '  -0.26%  '

# Row 45
Set-TextCell 'B45' 'Stellar'
Set-TextCell 'C45' 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextCell 'D45' '0.0986'
Set-TextCell 'E45' '  +0.98%  '

# Row 46
Set-TextCell 'B46' 'Mantle'
Set-TextCell 'C46' 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-TextCell 'D46' '0.610'
Set-TextCell 'E46' '  +2.31%  '

# Row 47
Set-TextCell 'D47' '19.78'
Set-TextCell 'E47' '  +3.14%  '

# Row 48
Set-TextCell 'D48' '0.0545'
Set-TextCell 'E48' '  +2.19%  '

# Row 49
Set-TextCell 'E49' '  +8.59%  '

# Row 50
Set-TextCell 'E50' '  +1.87%  '

# Row 51
Set-TextCell 'D51' '10.71'
Set-TextCell 'E51' '  +0.70%  '
